# Auto-generated edit script: updates cryptos list price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values; force text to avoid Excel numeric auto-conversion ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.513.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.836.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.662"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.741"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000318"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.432.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.831.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.129"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.287.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "93.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "47.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0988"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "638.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.429"
$ws.Range("D38").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0466"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.143"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.879.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000274"
$ws.Range("D51").Style = "Normal"

# --- Update Volume(1h) (column E) values ---
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("E3").Value = "  -4.08%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("E7").Value = "  -3.45%  "
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("E18").Value = "  -5.85%  "
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("E24").Value = "  -4.71%  "
$ws.Range("E25").Value = "  -3.81%  "
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  -9.26%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  -4.44%  "
$ws.Range("E31").Value = "  +5.56%  "
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("E35").Value = "  +8.02%  "
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("E37").Value = "  -4.80%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E43").Value = "  +11.60%  "
$ws.Range("E44").Value = "  -4.11%  "
$ws.Range("E45").Value = "  -5.74%  "
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("E47").Value = "  -4.79%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -17.04%  "
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("E51").Value = "  +0.47%  "

# --- Row 41/42: FirstDigitalUSD and dogwifhat swapped positions with updated figures ---
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +24.12%  "

